{"js": "// The title paragraph currently reads two runs:\n//   \"KRISTIAN BECKMAN AWARD \" (trailing space) + \"<year>\" (yellow highlight)\n// The edit removes the \"<year>\" placeholder run entirely and drops the\n// trailing space left behind, leaving a single run that reads exactly\n// \"KRISTIAN BECKMAN AWARD\".\nconst body = context.document.body;\n\n// Match the space + placeholder text uniquely (the other, unrelated\n// \"KRISTIAN BECKMAN AWARD \" occurrence earlier in the doc has no \"<year>\"\n// after it, so searching for \" <year>\" pins us to the right paragraph).\nconst hits = body.search(\" <year>\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  // Deleting the whole \" <year>\" span both removes the placeholder run and\n  // trims the trailing space off the preceding \"KRISTIAN BECKMAN AWARD\" run.\n  hits.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# The title paragraph currently reads two runs:\n#   \"KRISTIAN BECKMAN AWARD \" (trailing space) + \"<year>\" (yellow highlight)\n# The edit removes the \"<year>\" placeholder run entirely and drops the\n# trailing space left behind, leaving a single run that reads exactly\n# \"KRISTIAN BECKMAN AWARD\".\n$d = $word.ActiveDocument\n\n# Match the space + placeholder text uniquely (the other, unrelated\n# \"KRISTIAN BECKMAN AWARD \" occurrence earlier in the doc has no \"<year>\"\n# after it, so searching for \" <year>\" pins us to the right paragraph).\n$find = $d.Content.Find\n$find.Text = \" <year>\"\n$find.Replacement.Text = \"\"\n\n# wdFindContinue = 1, wdReplaceAll = 2.\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
